$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1194
$ws.Range("I19").Value = 588
$ws.Range("K19").Value = 588
$ws.Range("M19").Value = -413

$ws.Range("H70").Value = 2787.625
$ws.Range("I70").Value = 1525.5
$ws.Range("K70").Value = 4576.5
$ws.Range("M70").Value = -4306.5

$ws.Range("H73").Value = 2787.625
$ws.Range("I73").Value = 1525.5
$ws.Range("K73").Value = 4576.5
$ws.Range("M73").Value = -3640.5

$ws.Range("H106").Value = 3760
$ws.Range("I106").Value = 3760
$ws.Range("K106").Value = 3760
$ws.Range("M106").Value = -3129

$ws.Range("H107").Value = 3382.0476
$ws.Range("J107").Value = 2187.75
$ws.Range("L107").Value = 2187.75
$ws.Range("N107").Value = -6027.75

$ws.Range("H112").Value = 1779.4286
$ws.Range("I112").Value = 2183.7144
$ws.Range("J112").Value = 1644.6666
$ws.Range("K112").Value = 6551.1432
$ws.Range("L112").Value = 4933.9998
$ws.Range("M112").Value = -5443.1432
$ws.Range("N112").Value = -7149.9998

$ws.Range("H132").Value = 10604.857
$ws.Range("I132").Value = 11457.685
$ws.Range("K132").Value = 34373.055
$ws.Range("M132").Value = -31843.055

$ws.Range("H135").Value = 5363.087
$ws.Range("I135").Value = 1096.5
$ws.Range("K135").Value = 9868.5
$ws.Range("M135").Value = -7333.5

$ws.Range("H137").Value = 4144.615
$ws.Range("I137").Value = 1741.2858
$ws.Range("J137").Value = 10262.182
$ws.Range("K137").Value = 5223.857400000001
$ws.Range("L137").Value = 30786.546
$ws.Range("M137").Value = -2673.857400000001
$ws.Range("N137").Value = -35886.546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 151268.84
$ws.Range("I32").Value = 161309.38
$ws.Range("K32").Value = 161309.38
$ws.Range("M32").Value = -161022.38

$ws.Range("H61").Value = 7607.579
$ws.Range("I61").Value = 8032.294
$ws.Range("J61").Value = 3997.5
$ws.Range("K61").Value = 8032.294
$ws.Range("L61").Value = 3997.5
$ws.Range("M61").Value = -7820.294
$ws.Range("N61").Value = -4421.5

$ws.Range("H132").Value = 4489.95
$ws.Range("I132").Value = 2949
$ws.Range("K132").Value = 8847
$ws.Range("M132").Value = -6317

$ws.Range("H136").Value = 7607.579
$ws.Range("I136").Value = 8032.294
$ws.Range("J136").Value = 3997.5
$ws.Range("K136").Value = 24096.882
$ws.Range("L136").Value = 11992.5
$ws.Range("M136").Value = -21546.882
$ws.Range("N136").Value = -17092.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4181.364
$ws.Range("I94").Value = 4527.7334
$ws.Range("J94").Value = 3439.1428
$ws.Range("K94").Value = 4527.7334
$ws.Range("L94").Value = 3439.1428
$ws.Range("M94").Value = -4076.7334
$ws.Range("N94").Value = -4341.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 145.52942
$ws.Range("I7").Value = 93.5
$ws.Range("K7").Value = 93.5
$ws.Range("M7").Value = 19.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 7774.625
$ws.Range("I3").Value = 3266
$ws.Range("J3").Value = 13571.429
$ws.Range("K3").Value = 9798
$ws.Range("L3").Value = 40714.287
$ws.Range("M3").Value = -9686
$ws.Range("N3").Value = -40938.287

$ws.Range("H4").Value = 4031760.2
$ws.Range("I4").Value = 10538783
$ws.Range("K4").Value = 31616349
$ws.Range("M4").Value = -31616237

$ws.Range("H6").Value = 1856.6
$ws.Range("I6").Value = 947
$ws.Range("K6").Value = 2841
$ws.Range("M6").Value = -2728

$ws.Range("H7").Value = 71
$ws.Range("I7").Value = 75.40000000000001
$ws.Range("J7").Value = 65.5
$ws.Range("K7").Value = 226.2
$ws.Range("L7").Value = 196.5
$ws.Range("M7").Value = -114.2
$ws.Range("N7").Value = -420.5

$ws.Range("H10").Value = 1257.625
$ws.Range("I10").Value = 11
$ws.Range("J10").Value = 3335.3333
$ws.Range("K10").Value = 33
$ws.Range("L10").Value = 10005.9999
$ws.Range("M10").Value = 106
$ws.Range("N10").Value = -10283.9999

$ws.Range("H13").Value = 312.5
$ws.Range("J13").Value = 312.5
$ws.Range("L13").Value = 937.5
$ws.Range("N13").Value = -1273.5

$ws.Range("H15").Value = 622.4
$ws.Range("I15").Value = 528.25
$ws.Range("J15").Value = 999
$ws.Range("K15").Value = 1584.75
$ws.Range("L15").Value = 2997
$ws.Range("M15").Value = -1444.75
$ws.Range("N15").Value = -3277

$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

$ws.Range("H17").Value = 5009.6665
$ws.Range("I17").Value = 10000
$ws.Range("J17").Value = 2514.5
$ws.Range("K17").Value = 30000
$ws.Range("L17").Value = 7543.5
$ws.Range("M17").Value = -29831
$ws.Range("N17").Value = -7881.5

$ws.Range("H21").Value = 1210.6
$ws.Range("J21").Value = 1701.3334
$ws.Range("L21").Value = 5104.0002
$ws.Range("N21").Value = -5450.0002

$ws.Range("H25").Value = 989.4
$ws.Range("J25").Value = 1422.5
$ws.Range("L25").Value = 4267.5
$ws.Range("N25").Value = -4605.5

$ws.Range("H26").Value = 198.61539
$ws.Range("I26").Value = 209.5
$ws.Range("J26").Value = 181.2
$ws.Range("K26").Value = 628.5
$ws.Range("L26").Value = 543.5999999999999
$ws.Range("M26").Value = -340.5
$ws.Range("N26").Value = -1119.6

$ws.Range("H29").Value = 375.5
$ws.Range("J29").Value = 250
$ws.Range("L29").Value = 750
$ws.Range("N29").Value = -1304

$ws.Range("H30").Value = 989.4
$ws.Range("J30").Value = 1422.5
$ws.Range("L30").Value = 4267.5
$ws.Range("N30").Value = -4471.5

$ws.Range("I38").Value = 700
$ws.Range("J38").Value = 165
$ws.Range("K38").Value = 2100
$ws.Range("L38").Value = 495
$ws.Range("M38").Value = -1753
$ws.Range("N38").Value = -1189

$ws.Range("H107").Value = 3654.8667
$ws.Range("J107").Value = 4418.636
$ws.Range("L107").Value = 13255.908
$ws.Range("N107").Value = -17095.908

$ws.Range("H137").Value = 7488.625
$ws.Range("I137").Value = 1612.5714
$ws.Range("J137").Value = 12058.889
$ws.Range("K137").Value = 4837.7142
$ws.Range("L137").Value = 36176.667
$ws.Range("M137").Value = 262.2857999999997
$ws.Range("N137").Value = -46376.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws.Range("H132").Value = 9566.861999999999
$ws.Range("I132").Value = 7472.35
$ws.Range("J132").Value = 14221.333
$ws.Range("K132").Value = 22417.05
$ws.Range("L132").Value = 42663.999
$ws.Range("M132").Value = -19887.05
$ws.Range("N132").Value = -47723.999

$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1198.8334
$ws.Range("I16").Value = 1258.8
$ws.Range("J16").Value = 899
$ws.Range("K16").Value = 1258.8
$ws.Range("L16").Value = 899
$ws.Range("M16").Value = -1088.8
$ws.Range("N16").Value = -1239

$ws.Range("H22").Value = 2879.0417
$ws.Range("I22").Value = 830.3333
$ws.Range("J22").Value = 3561.9443
$ws.Range("K22").Value = 830.3333
$ws.Range("L22").Value = 3561.9443
$ws.Range("M22").Value = -535.3333
$ws.Range("N22").Value = -4151.9443

$ws.Range("H27").Value = 2879.0417
$ws.Range("I27").Value = 830.3333
$ws.Range("J27").Value = 3561.9443
$ws.Range("K27").Value = 830.3333
$ws.Range("L27").Value = 3561.9443
$ws.Range("M27").Value = -723.3333
$ws.Range("N27").Value = -3775.9443

$ws.Range("H122").Value = 5867.6665
$ws.Range("I122").Value = 3286.3333
$ws.Range("J122").Value = 7158.3335
$ws.Range("K122").Value = 9858.999899999999
$ws.Range("L122").Value = 21475.0005
$ws.Range("M122").Value = -7408.999899999999
$ws.Range("N122").Value = -26375.0005

$ws.Range("H132").Value = 3821.111
$ws.Range("I132").Value = 4250
$ws.Range("K132").Value = 12750
$ws.Range("M132").Value = -10220

$ws.Range("H136").Value = 7749.8335
$ws.Range("I136").Value = 5500
$ws.Range("K136").Value = 16500
$ws.Range("M136").Value = -13950

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 49197
$ws.Range("I82").Value = 26000
$ws.Range("J82").Value = 54996.25
$ws.Range("K82").Value = 26000
$ws.Range("L82").Value = 54996.25
$ws.Range("M82").Value = -25617
$ws.Range("N82").Value = -55762.25

$ws.Range("H85").Value = 49197
$ws.Range("I85").Value = 26000
$ws.Range("J85").Value = 54996.25
$ws.Range("K85").Value = 26000
$ws.Range("L85").Value = 54996.25
$ws.Range("M85").Value = -24674
$ws.Range("N85").Value = -57648.25

$ws.Range("H104").Value = 19975
$ws.Range("J104").Value = 19975
$ws.Range("L104").Value = 19975
$ws.Range("N104").Value = -26963

$ws.Range("H113").Value = 2085.2285
$ws.Range("I113").Value = 1622
$ws.Range("K113").Value = 4866
$ws.Range("M113").Value = -2696

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H126").Value = 2269.889
$ws.Range("I126").Value = 1945
$ws.Range("K126").Value = 5835
$ws.Range("M126").Value = -3365

$ws.Range("H132").Value = 2299.7036
$ws.Range("I132").Value = 2293.8696
$ws.Range("K132").Value = 6881.6088
$ws.Range("M132").Value = -4351.6088

$ws.Range("H136").Value = 127805.625
$ws.Range("I136").Value = 2906.6667
$ws.Range("K136").Value = 8720.000100000001
$ws.Range("M136").Value = -6170.000100000001
